$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.255.61'
$ws.Range('E2').Value = '  -3.54%  '
$ws.Range('D3').Value = '2.466.70'
$ws.Range('E3').Value = '  -2.81%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '311.70'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').Value = '94.44'
$ws.Range('E6').Value = '  -6.13%  '
$ws.Range('E7').Value = '  -3.35%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.500'
$ws.Range('E9').Value = '  -4.45%  '
$ws.Range('D10').Value = '33.51'
$ws.Range('E10').Value = '  -6.07%  '
$ws.Range('D11').Value = '0.0782'
$ws.Range('E11').Value = '  -2.83%  '
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').Value = '6.99'
$ws.Range('E13').Value = '  -4.56%  '
$ws.Range('D14').Value = '2.841.68'
$ws.Range('D15').Value = '2.476.64'
$ws.Range('E15').Value = '  -3.18%  '
$ws.Range('D16').Value = '15.03'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('D17').Value = '0.788'
$ws.Range('E17').Value = '  -3.25%  '
$ws.Range('D18').Value = '41.238.03'
$ws.Range('E18').Value = '  -3.59%  '
$ws.Range('D19').Value = '6.32'
$ws.Range('E19').Value = '  -6.37%  '
$ws.Range('D20').Value = '0.0₃0924'
$ws.Range('E20').Value = '  -2.98%  '
$ws.Range('E21').Value = '  -8.98%  '
$ws.Range('D22').Value = '68.42'
$ws.Range('E22').Value = '  -2.46%  '
$ws.Range('D23').Value = '237.29'
$ws.Range('E24').Value = '  -4.86%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  -6.48%  '
$ws.Range('D27').Value = '24.01'
$ws.Range('E27').Value = '  -5.73%  '
$ws.Range('E28').Value = '  -4.29%  '
$ws.Range('D29').Value = '9.67'
$ws.Range('E29').Value = '  -4.79%  '
$ws.Range('E30').Value = '  -5.88%  '
$ws.Range('D31').Value = '151.53'
$ws.Range('E31').Value = '  -4.69%  '
$ws.Range('D32').Value = '5.47'
$ws.Range('E32').Value = '  -6.51%  '
$ws.Range('D33').Value = '2.60'
$ws.Range('E33').Value = '  -5.62%  '
$ws.Range('E34').Value = '  -3.15%  '
$ws.Range('D35').Value = '0.0745'
$ws.Range('E35').Value = '  -5.99%  '
$ws.Range('E36').Value = '  -3.16%  '
$ws.Range('D37').Value = '17.29'
$ws.Range('E37').Value = '  -4.61%  '
$ws.Range('D38').Value = '1.87'
$ws.Range('E38').Value = '  -5.02%  '
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '4.26'
$ws.Range('E40').Value = '  +3.00%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.101'
$ws.Range('E41').Value = '  -8.89%  '
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = '19.41'
$ws.Range('E43').Value = '  -10.88%  '
$ws.Range('D44').Value = '1.979.61'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('E45').Value = '  -4.60%  '
$ws.Range('E46').Value = '  -8.54%  '
$ws.Range('D47').Value = '8.67'
$ws.Range('E47').Value = '  -5.52%  '
$ws.Range('D48').Value = '2.705.17'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('E49').Value = '  -3.88%  '
$ws.Range('D50').Value = '96.65'
$ws.Range('E50').Value = '  -4.53%  '
$ws.Range('D51').Value = '74.56'
$ws.Range('E51').Value = '  -6.85%  '
